$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 560.70966
$ws.Range("I28").Value = 596.5909
$ws.Range("J28").Value = 473
$ws.Range("K28").Value = 596.5909
$ws.Range("L28").Value = 473
$ws.Range("M28").Value = -111.5909
$ws.Range("N28").Value = -1443
$ws.Range("H98").Value = 2267.3157
$ws.Range("I98").Value = 879
$ws.Range("J98").Value = 5275.3335
$ws.Range("K98").Value = 879
$ws.Range("L98").Value = 5275.3335
$ws.Range("M98").Value = 619
$ws.Range("N98").Value = -8271.333500000001
$ws.Range("H107").Value = 1705.6364
$ws.Range("I107").Value = 1433.7646
$ws.Range("J107").Value = 2630
$ws.Range("K107").Value = 1433.7646
$ws.Range("L107").Value = 2630
$ws.Range("M107").Value = 486.2354
$ws.Range("N107").Value = -6470
$ws.Range("H112").Value = 1507.8462
$ws.Range("J112").Value = 1588.8914
$ws.Range("L112").Value = 4766.674199999999
$ws.Range("N112").Value = -6982.674199999999
$ws.Range("H121").Value = 2319.7778
$ws.Range("J121").Value = 2319.7778
$ws.Range("L121").Value = 6959.3334
$ws.Range("N121").Value = -10453.3334
$ws.Range("H122").Value = 2267.3157
$ws.Range("I122").Value = 879
$ws.Range("J122").Value = 5275.3335
$ws.Range("K122").Value = 2637
$ws.Range("L122").Value = 15826.0005
$ws.Range("M122").Value = -187
$ws.Range("N122").Value = -20726.0005
$ws.Range("H138").Value = 3354.847
$ws.Range("I138").Value = 743.7560999999999
$ws.Range("J138").Value = 5233
$ws.Range("K138").Value = 2231.2683
$ws.Range("L138").Value = 15699
$ws.Range("M138").Value = 2908.7317
$ws.Range("N138").Value = -25979
$ws.Range("H141").Value = 26328.432
$ws.Range("I141").Value = 28576.15
$ws.Range("K141").Value = 85728.45000000001
$ws.Range("M141").Value = -80548.45000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 9747.5
$ws.Range("J3").Value = 9747.5
$ws.Range("L3").Value = 9747.5
$ws.Range("N3").Value = -9977.5
$ws.Range("H74").Value = 3488.2646
$ws.Range("I74").Value = 3665.5
$ws.Range("K74").Value = 3665.5
$ws.Range("M74").Value = -2791.5
$ws.Range("H77").Value = 3488.2646
$ws.Range("I77").Value = 3665.5
$ws.Range("K77").Value = 18327.5
$ws.Range("M77").Value = -13959.5
$ws.Range("H132").Value = 1525.5272
$ws.Range("I132").Value = 895.7778
$ws.Range("K132").Value = 2687.3334
$ws.Range("M132").Value = -157.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10001966
$ws.Range("I31").Value = 1204.9678
$ws.Range("J31").Value = 26318998
$ws.Range("K31").Value = 1204.9678
$ws.Range("L31").Value = 26318998
$ws.Range("M31").Value = -909.9677999999999
$ws.Range("N31").Value = -26319588
$ws.Range("H34").Value = 10001966
$ws.Range("I34").Value = 1204.9678
$ws.Range("J34").Value = 26318998
$ws.Range("K34").Value = 1204.9678
$ws.Range("L34").Value = 26318998
$ws.Range("M34").Value = -1002.9678
$ws.Range("N34").Value = -26319402
$ws.Range("H99").Value = 11116177
$ws.Range("I99").Value = 22225076
$ws.Range("J99").Value = 7277.778
$ws.Range("K99").Value = 22225076
$ws.Range("L99").Value = 7277.778
$ws.Range("M99").Value = -22223578
$ws.Range("N99").Value = -10273.778
$ws.Range("H126").Value = 11116177
$ws.Range("I126").Value = 22225076
$ws.Range("J126").Value = 7277.778
$ws.Range("K126").Value = 66675228
$ws.Range("L126").Value = 21833.334
$ws.Range("M126").Value = -66672758
$ws.Range("N126").Value = -26773.334
$ws.Range("H132").Value = 5960.722
$ws.Range("I132").Value = 5869.5
$ws.Range("J132").Value = 6074.75
$ws.Range("K132").Value = 17608.5
$ws.Range("L132").Value = 18224.25
$ws.Range("M132").Value = -15078.5
$ws.Range("N132").Value = -23284.25
$ws.Range("H134").Value = 7790.0527
$ws.Range("I134").Value = 10951.1
$ws.Range("J134").Value = 4277.778
$ws.Range("K134").Value = 32853.3
$ws.Range("L134").Value = 12833.334
$ws.Range("M134").Value = -30318.3
$ws.Range("N134").Value = -17903.334
$ws.Range("H141").Value = 33938.89
$ws.Range("J141").Value = 33938.89
$ws.Range("L141").Value = 33938.89
$ws.Range("N141").Value = -44298.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 29870.03
$ws.Range("I107").Value = 465.8
$ws.Range("J107").Value = 53083.895
$ws.Range("K107").Value = 1397.4
$ws.Range("L107").Value = 159251.685
$ws.Range("M107").Value = 522.5999999999999
$ws.Range("N107").Value = -163091.685
$ws.Range("H137").Value = 2884
$ws.Range("J137").Value = 3946.923
$ws.Range("L137").Value = 11840.769
$ws.Range("N137").Value = -22040.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 28499.5
$ws.Range("J4").Value = 28499.5
$ws.Range("L4").Value = 28499.5
$ws.Range("N4").Value = -28723.5
$ws.Range("H113").Value = 3342.2
$ws.Range("I113").Value = 3342.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3342.2
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1172.2
$ws.Range("N113").ClearContents()
$ws.Range("H123").Value = 11674.5625
$ws.Range("J123").Value = 11674.5625
$ws.Range("L123").Value = 11674.5625
$ws.Range("N123").Value = -16574.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 12500
$ws.Range("J39").Value = 12500
$ws.Range("L39").Value = 12500
$ws.Range("N39").Value = -13420
$ws.Range("H42").Value = 39996.332
$ws.Range("J42").Value = 39996.332
$ws.Range("L42").Value = 39996.332
$ws.Range("N42").Value = -41122.332
$ws.Range("H49").Value = 39996.332
$ws.Range("J49").Value = 39996.332
$ws.Range("L49").Value = 39996.332
$ws.Range("N49").Value = -40290.332
$ws.Range("H123").Value = 31609.334
$ws.Range("J123").Value = 31609.334
$ws.Range("L123").Value = 31609.334
$ws.Range("N123").Value = -41409.334
$ws.Range("H132").Value = 3790.775
$ws.Range("I132").Value = 1685.8846
$ws.Range("J132").Value = 7699.857
$ws.Range("K132").Value = 5057.6538
$ws.Range("L132").Value = 23099.571
$ws.Range("M132").Value = -2527.6538
$ws.Range("N132").Value = -28159.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 29516.666
$ws.Range("J92").Value = 29516.666
$ws.Range("L92").Value = 29516.666
$ws.Range("N92").Value = -34508.666
$ws.Range("H136").Value = 7105.5
$ws.Range("I136").Value = 3517
$ws.Range("J136").Value = 8301.666999999999
$ws.Range("K136").Value = 10551
$ws.Range("L136").Value = 24905.001
$ws.Range("M136").Value = -8001
$ws.Range("N136").Value = -30005.001
$ws.Range("H138").Value = 43222
$ws.Range("J138").Value = 43222
$ws.Range("L138").Value = 43222
$ws.Range("N138").Value = -53502
$ws.Range("H139").Value = 37140.535
$ws.Range("J139").Value = 37010.555
$ws.Range("L139").Value = 37010.555
$ws.Range("N139").Value = -47290.555
$ws.Range("H140").Value = 46575.285
$ws.Range("J140").Value = 46575.285
$ws.Range("L140").Value = 46575.285
$ws.Range("N140").Value = -56935.285
$ws.Range("H141").Value = 43142.855
$ws.Range("J141").Value = 43142.855
$ws.Range("L141").Value = 43142.855
$ws.Range("N141").Value = -53502.855
